# Accession Register page code....
#
# Adds a new "AccessionRegister" worksheet at the end of the workbook,
# populated with a header row and a single sample data row, mirroring the
# BookCatalog-style sheets already present in this workbook.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet as the last tab -----------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "AccessionRegister"

# --- Header row (row 1) -----------------------------------------------------
$headers = @(
    "BookTitle",
    "AccessionNo",
    "Author",
    "Place",
    "Publisher",
    "Year",
    "Pages",
    "Sources",
    "BillNo",
    "Cost",
    "ClassNo",
    "BookNo",
    "WithdrawNo"
)
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- Sample data row (row 2) ------------------------------------------------
$ws.Cells.Item(2, 1).Value = "Automation"
$ws.Cells.Item(2, 2).Value = "Ac232"
$ws.Cells.Item(2, 3).Value = "Steven"
$ws.Cells.Item(2, 4).Value = "Europe"
$ws.Cells.Item(2, 5).Value = "Wills"
$ws.Cells.Item(2, 6).Value = 2018
$ws.Cells.Item(2, 7).Value = 200
$ws.Cells.Item(2, 8).Value = "Printed"
$ws.Cells.Item(2, 9).Value = "BL30"
$ws.Cells.Item(2, 10).Value = 400
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = "BK303"
$ws.Cells.Item(2, 13).Value = "WD20"

# --- Column widths (best match obtainable through ColumnWidth) -------------
$widths = @(15.67, 22.33, 16.0, 18.67, 22.83, 19.33, 18.33, 15.67, 21.33, 17.0, 18.67, 21.17, 19.67)
for ($i = 0; $i -lt $widths.Count; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i]
}

# --- View state: new sheet becomes the active / selected tab ---------------
$ws.Range("M2").Select() | Out-Null
